$wb = $excel.ActiveWorkbook

# --- Sheet2 (users) data edits ---
$ws2 = $wb.Worksheets.Item("Sheet2")

# Row 2 -> steve's data (previously row 3)
$ws2.Range("A2").Value = "steve"
$ws2.Range("B2").Value = "steve@123.com"
$ws2.Range("C2").Value = "'1234567890"
$ws2.Range("D2").Value = "'25"
$ws2.Range("E2").Value = "steve@123"
$ws2.Range("C2:D2").ClearFormats()

# Row 3 -> mani's data (modified)
$ws2.Range("A3").Value = "mani"
$ws2.Range("B3").Value = "mani@123.com"
$ws2.Range("C3").Value = "'987654321"
$ws2.Range("D3").Value = "'20"
$ws2.Range("E3").Value = "mani123"
$ws2.Range("C3:D3").ClearFormats()

# Row 4 (sai) removed entirely
$ws2.Rows.Item(4).Delete()

# Select A3 on Sheet2 and make it the active/tabSelected sheet
$ws2.Range("A3").Select()
$ws2.Activate()

$wb.Save()
